$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D2 (Instance): Automation3 -> Automation2
$ws.Range("D2").Value = "Automation2"

# B2 (TestCases): 34 -> 42,43,44,45
$ws.Range("B2").Value = "42,43,44,45"
